$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.9
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 4.33
$ws.Range("J3").Value = 2.63
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("Z3").Value = 15
$ws.Range("AH3").Value = 21
$ws.Range("AJ3").Value = 51
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 151
$ws.Range("G5").Value = 2.05
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 2.88
$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 4.75
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("X5").Value = 8.5
$ws.Range("Y5").Value = 9.5
$ws.Range("Z5").Value = 17
$ws.Range("AA5").Value = 21
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 6.5
$ws.Range("AD5").Value = 6
$ws.Range("AE5").Value = 19
$ws.Range("AI5").Value = 15
$ws.Range("AK5").Value = 41
$ws.Range("AL5").Value = 51
$ws.Range("AM5").Value = 900
$ws.Range("AO5").Value = 12
$ws.Range("AP5").Value = 26
$ws.Range("AQ5").Value = 41
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AW5").Value = 5.5
$ws.Range("AX5").Value = 23
$ws.Range("BB5").Value = 351
$ws.Range("N7").Value = 13.8
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.62
$ws.Range("R9").Value = 2.25
$ws.Range("G15").Value = 2.55
$ws.Range("M15").Value = 1.08
$ws.Range("N15").Value = 8
$ws.Range("Q15").Value = 2.3
$ws.Range("R15").Value = 1.6
$ws.Range("U15").Value = 1.91
$ws.Range("V15").Value = 1.8
$ws.Range("AB15").Value = 34
$ws.Range("AE15").Value = 15
$ws.Range("AK15").Value = 26
$ws.Range("AM15").Value = 800
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 9
$ws.Range("Q16").Value = 2.1
$ws.Range("R16").Value = 1.7
$ws.Range("O17").Value = 1.18
$ws.Range("P17").Value = 4.5
$ws.Range("Q17").Value = 1.65
$ws.Range("R17").Value = 2.2
$ws.Range("G18").Value = 2.25
$ws.Range("H18").Value = 2.9
$ws.Range("I18").Value = 3.7
$ws.Range("J18").Value = 3.1
$ws.Range("L18").Value = 4.33
$ws.Range("M18").Value = 1.13
$ws.Range("N18").Value = 6
$ws.Range("O18").Value = 1.57
$ws.Range("P18").Value = 2.25
$ws.Range("Q18").Value = 2.88
$ws.Range("R18").Value = 1.4
$ws.Range("S18").Value = 1.62
$ws.Range("T18").Value = 2.2
$ws.Range("U18").Value = 2.25
$ws.Range("V18").Value = 1.57
$ws.Range("W18").Value = 5.5
$ws.Range("Z18").Value = 21
$ws.Range("AA18").Value = 23
$ws.Range("AC18").Value = 5.5
$ws.Range("AE18").Value = 21
$ws.Range("AF18").Value = 81
$ws.Range("AG18").Value = 8
$ws.Range("AK18").Value = 41
$ws.Range("AQ18").Value = 51
$ws.Range("AS18").Value = 301
$ws.Range("AT18").Value = 2.2
$ws.Range("AU18").Value = 9.5
$ws.Range("AY18").Value = 41
$ws.Range("G19").Value = 2.05
$ws.Range("H19").Value = 3.3
$ws.Range("I19").Value = 3.6
$ws.Range("J19").Value = 2.75
$ws.Range("K19").Value = 2.1
$ws.Range("L19").Value = 4
$ws.Range("M19").Value = 1.06
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 1.3
$ws.Range("P19").Value = 3.4
$ws.Range("Q19").Value = 2.03
$ws.Range("R19").Value = 1.83
$ws.Range("X19").Value = 10
$ws.Range("Z19").Value = 19
$ws.Range("AA19").Value = 17
$ws.Range("AB19").Value = 29
$ws.Range("AC19").Value = 9.5
$ws.Range("AG19").Value = 10
$ws.Range("AH19").Value = 17
$ws.Range("AK19").Value = 29
$ws.Range("AL19").Value = 34
$ws.Range("AO19").Value = 11
$ws.Range("AQ19").Value = 41
$ws.Range("AX19").Value = 19
$ws.Range("G20").Value = 2.25
$ws.Range("H20").Value = 2.9
$ws.Range("I20").Value = 3.7
$ws.Range("J20").Value = 3
$ws.Range("L20").Value = 4
$ws.Range("W20").Value = 7
$ws.Range("X20").Value = 10
$ws.Range("Y20").Value = 9.5
$ws.Range("Z20").Value = 21
$ws.Range("AA20").Value = 19
$ws.Range("AG20").Value = 9.5
$ws.Range("AH20").Value = 17
$ws.Range("AI20").Value = 13
$ws.Range("AJ20").Value = 41
$ws.Range("AN20").Value = 4
$ws.Range("AQ20").Value = 41
$ws.Range("AX20").Value = 19
$ws.Range("AZ20").Value = 67
$ws.Range("BA20").Value = 101
$ws.Range("AM21").Value = 1000
$ws.Range("W22").Value = 13
$ws.Range("Z22").Value = 67
$ws.Range("AF22").Value = 81
$ws.Range("G25").Value = 2.05
$ws.Range("I25").Value = 3.1
$ws.Range("J25").Value = 2.65
$ws.Range("L25").Value = 3.6
$ws.Range("O25").Value = 1.25
$ws.Range("P25").Value = 3.6
$ws.Range("R25").Value = 2
$ws.Range("U25").Value = 1.65
$ws.Range("V25").Value = 2.12
$ws.Range("X25").Value = 10.75
$ws.Range("Z25").Value = 19.5
$ws.Range("AA25").Value = 15.5
$ws.Range("AF25").Value = 55
$ws.Range("AG25").Value = 11
$ws.Range("AI25").Value = 11
$ws.Range("AJ25").Value = 40
$ws.Range("AK25").Value = 25
$ws.Range("AL25").Value = 30
$ws.Range("AN25").Value = 4.1
$ws.Range("AO25").Value = 10.5
$ws.Range("AP25").Value = 17.5
$ws.Range("AR25").Value = 65
$ws.Range("AU25").Value = 6.9
$ws.Range("AW25").Value = 5.2
$ws.Range("AX25").Value = 16.5
$ws.Range("AY25").Value = 23
$ws.Range("AZ25").Value = 80
$ws.Range("BA25").Value = 110
$ws.Range("G27").Value = 2.22
$ws.Range("H27").Value = 2.82
$ws.Range("I27").Value = 3.5
$ws.Range("J27").Value = 2.87
$ws.Range("L27").Value = 4.2
$ws.Range("O27").Value = 1.5
$ws.Range("P27").Value = 2.27
$ws.Range("Q27").Value = 2.4
$ws.Range("R27").Value = 1.44
$ws.Range("S27").Value = 1.53
$ws.Range("T27").Value = 2.2
$ws.Range("W27").Value = 5.9
$ws.Range("X27").Value = 9.75
$ws.Range("Y27").Value = 9.25
$ws.Range("Z27").Value = 22
$ws.Range("AA27").Value = 21
$ws.Range("AB27").Value = 40
$ws.Range("AC27").Value = 6.2
$ws.Range("AG27").Value = 7.5
$ws.Range("AH27").Value = 17
$ws.Range("AI27").Value = 13
$ws.Range("AJ27").Value = 55
$ws.Range("AK27").Value = 40
$ws.Range("AL27").Value = 60
$ws.Range("AN27").Value = 3.9
$ws.Range("AO27").Value = 12
$ws.Range("AP27").Value = 23
$ws.Range("AQ27").Value = 50
$ws.Range("AR27").Value = 100
$ws.Range("AT27").Value = 2.18
$ws.Range("AW27").Value = 5.2
$ws.Range("AX27").Value = 22
$ws.Range("AY27").Value = 32
$ws.Range("AZ27").Value = 120
